# Links.xlsx edit — add HoverDroids Spin Library links alongside the existing
# Parallax Spin Library links, plus a header row ("Title"/"Slug"/link columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new header row at the top; everything else shifts down ---
$ws.Rows.Item(1).Insert()

$ws.Range("A1").Value = "Title"
$ws.Range("B1").Value = "Slug"
$ws.Range("C1").Value = "Link (Parallax Spin Library)"
$ws.Range("D1").Value = "Link (HoverDroids Spin Library)"

# --- 2. Re-create the Parallax hyperlink, now on E2 (it did not move itself
#        when the row above it was inserted) ---
$ws.Range("E1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("E2"), "http://propeller-microcontroller.hoverdroids.com/parallax-propeller-spin-library/") | Out-Null

# --- 3. Add the new HoverDroids base-url cell + hyperlink on E3 ---
$ws.Range("E3").Value = "http://propeller-microcontroller.hoverdroids.com/propeller-spin-library/"
$ws.Hyperlinks.Add($ws.Range("E3"), "http://propeller-microcontroller.hoverdroids.com/propeller-spin-library/") | Out-Null

# --- 4. Build column D: same <li><a href='...'>Title</a></li> formula as
#        column C, but pointing at the HoverDroids base url in $E$3 ---
$ws.Range("D2").Formula = "=""<li><a href='"" & `$E`$3 & B2 & ""'>"" &A2&""</a></li>"""
$ws.Range("D3:D56").Formula = "=""<li><a href='"" & `$E`$3 & B3 & ""'>"" &A3&""</a></li>"""

# --- 5. Clear the stray direct formatting (a leftover font / wrap-text style)
#        that a few cells carried in the original sheet ---
$clearCells = @("B4","B16","A17","B17","B22","B31","B32","B44","B51","B52","B53","B54","B55","A56","B56")
foreach ($addr in $clearCells) {
    $ws.Range($addr).ClearFormats()
}

# --- 6. Column widths / best-fit ---
$ws.Columns.Item(1).ColumnWidth = 33.5
$ws.Columns.Item(3).ColumnWidth = 154.5
$ws.Columns.Item(4).ColumnWidth = 111.5

# --- 7. View state ---
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Application.ActiveWindow.ScrollRow = 40
$ws.Range("D57").Select()

Write-Output "edit applied"
